$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: "Parent Phone" -> "parent phone"
$ws.Range("E1").Value = "parent phone"

# Body values: 0 -> 123 for all data rows (E2:E59)
$ws.Range("E2:E59").Value = 123

# Update the view/selection to match the saved workbook state
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("E2:E59").Select()
